$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume %), with a few coin rows
# reordered to reflect the updated ranking.

$ws.Range("D2").Value = "55.083.32"
$ws.Range("E2").Value = "  +7.31%  "

$ws.Range("D3").Value = "2.458.04"
$ws.Range("E3").Value = "  +8.92%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "478.49"
$ws.Range("E5").Value = "  +10.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.10"
$ws.Range("E6").Value = "  +18.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("E8").Value = "  +9.88%  "

$ws.Range("D9").Value = "2.457.75"
$ws.Range("E9").Value = "  +8.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0960"
$ws.Range("E10").Value = "  +13.97%  "

$ws.Range("E11").Value = "  +6.16%  "

$ws.Range("E12").Value = "  +9.05%  "

$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").Value = "2.862.17"
$ws.Range("E14").Value = "  +8.71%  "

$ws.Range("D15").Value = "55.155.79"
$ws.Range("E15").Value = "  +7.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.42"
$ws.Range("E16").Value = "  +10.93%  "

$ws.Range("E17").Value = "  +18.18%  "

$ws.Range("D18").Value = "2.453.01"
$ws.Range("E18").Value = "  +9.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.34"
$ws.Range("E19").Value = "  +11.90%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.92"
$ws.Range("E20").Value = "  +16.30%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.45"
$ws.Range("E21").Value = "  +6.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.994"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("E23").Value = "  +12.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.16"
$ws.Range("E24").Value = "  +8.86%  "

$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("E26").Value = "  +11.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +18.58%  "

$ws.Range("D28").Value = "2.555.91"
$ws.Range("E28").Value = "  +8.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  +9.20%  "

$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  +23.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.70"
$ws.Range("E32").Value = "  +4.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.88"
$ws.Range("E33").Value = "  +9.25%  "

$ws.Range("E34").Value = "  +13.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.14"
$ws.Range("E35").Value = "  +12.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("E36").Value = "  +14.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.847"
$ws.Range("E37").Value = "  +9.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  +8.84%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.33"
$ws.Range("E40").Value = "  +5.66%  "

$ws.Range("E41").Value = "  +9.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  +11.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0540"
$ws.Range("E43").Value = "  +10.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.28"
$ws.Range("E44").Value = "  +13.12%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.12"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "255.91"
$ws.Range("E46").Value = "  +33.74%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.64"
$ws.Range("E47").Value = "  +16.99%  "

$ws.Range("E48").Value = "  +11.44%  "

$ws.Range("D49").Value = "1.926.81"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  +10.93%  "

$ws.Range("E51").Value = "  +11.64%  "
